$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44984
$ws.Range("J2").Value = 200
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17500
$ws.Range("P2").Value = 972

# Row 3
$ws.Range("D3").Value = 45068
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 16000
$ws.Range("M3").Value = 16500
$ws.Range("P3").Value = 917

# Row 4
$ws.Range("D4").Value = 44957
$ws.Range("K4").Value = 21000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 21500
$ws.Range("P4").Value = 1194

# Row 5
$ws.Range("D5").Value = 45230
$ws.Range("J5").Value = 360
$ws.Range("K5").Value = 16000
$ws.Range("L5").Value = 17000
$ws.Range("M5").Value = 16500
$ws.Range("P5").Value = 917

# Row 6
$ws.Range("D6").Value = 45117
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 17000
$ws.Range("L6").Value = 18000
$ws.Range("M6").Value = 17500
$ws.Range("P6").Value = 972

# Row 7
$ws.Range("D7").Value = 44557
$ws.Range("K7").Value = 13000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 13500
$ws.Range("P7").Value = 750

# Row 8
$ws.Range("D8").Value = 45152
$ws.Range("J8").Value = 500

# Row 9
$ws.Range("D9").Value = 44998
$ws.Range("J9").Value = 320
$ws.Range("K9").Value = 17000
$ws.Range("L9").Value = 18000
$ws.Range("M9").Value = 17500
$ws.Range("P9").Value = 972

# Row 10
$ws.Range("D10").Value = 45194

# Row 12
$ws.Range("D12").Value = 45222
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 16000
$ws.Range("L12").Value = 17000
$ws.Range("M12").Value = 16500
$ws.Range("P12").Value = 917

# Row 13
$ws.Range("D13").Value = 44977
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 16500
$ws.Range("L13").Value = 17000
$ws.Range("M13").Value = 16750
$ws.Range("P13").Value = 931

# Row 14
$ws.Range("D14").Value = 44964
$ws.Range("J14").Value = 300
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 21000
$ws.Range("M14").Value = 20500
$ws.Range("P14").Value = 1139

# Row 15
$ws.Range("D15").Value = 45005
$ws.Range("J15").Value = 200

# Row 16
$ws.Range("D16").Value = 45180

# Row 17
$ws.Range("D17").Value = 45154
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 16500
$ws.Range("M17").Value = 16750
$ws.Range("P17").Value = 931

# Row 18
$ws.Range("D18").Value = 45159
$ws.Range("J18").Value = 400

# Row 19
$ws.Range("D19").Value = 45229
$ws.Range("J19").Value = 460
$ws.Range("K19").Value = 16000
$ws.Range("L19").Value = 17000
$ws.Range("M19").Value = 16500
$ws.Range("P19").Value = 917

# Row 20
$ws.Range("D20").Value = 44960
$ws.Range("J20").Value = 400
$ws.Range("K20").Value = 19500
$ws.Range("L20").Value = 20000
$ws.Range("M20").Value = 19750
$ws.Range("P20").Value = 1097

# Row 22
$ws.Range("D22").Value = 44568
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 16000
$ws.Range("M22").Value = 15500
$ws.Range("P22").Value = 861

# Row 23
$ws.Range("D23").Value = 44547
$ws.Range("J23").Value = 200
$ws.Range("K23").Value = 13000
$ws.Range("L23").Value = 14000
$ws.Range("M23").Value = 13500
$ws.Range("P23").Value = 750

# Row 24
$ws.Range("D24").Value = 45142
$ws.Range("K24").Value = 17000
$ws.Range("L24").Value = 18000
$ws.Range("M24").Value = 17500
$ws.Range("P24").Value = 972

# Row 25
$ws.Range("D25").Value = 45177
$ws.Range("J25").Value = 540
$ws.Range("K25").Value = 16000
$ws.Range("L25").Value = 17000
$ws.Range("M25").Value = 16500
$ws.Range("P25").Value = 917

